$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Ensure price/volume columns are treated as text so values like "1.000" or "0.9998" are preserved exactly
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.105.96"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.782.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.19%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.03"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4926"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2678"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06259"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.782.54"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07037"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6270"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.642"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "80.14"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.086.09"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9999"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007225"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.00"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.004.50"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.573"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.720"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.235"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.01"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.80"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.864"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "109.16"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.389"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.174"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +6.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08291"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.785"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04915"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +9.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.073"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6528"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +4.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.609"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9491"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +7.02%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.053"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.949"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.64%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01552"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.38%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.85"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3989"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.188"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1204"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05426"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.004"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.34%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.53"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.88"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.02%  "
